$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly hours log currently ends at row 14 with a SUM total in B14
# (SUM(B5:B13)). We're adding one more day of data ("Sunday", 0.75 hrs)
# as the new row 14, and pushing the running total down to row 15 so it
# now sums B5:B14.

# 1) Push the total row's formatting/formula down one row first, while
#    B14 still holds the SUM cell's style (s="4").
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)

# 2) Give the new data row (row 14) the same look as the other weekday
#    rows by copying formats from row 13 (day name style + hours style).
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

# 3) Fill in the new day's data.
$ws.Range("A14").Value = "Sunday"
$ws.Range("B14").Value = 0.75

# 4) Extend the running total so it includes the new row, now in B15.
$ws.Range("B15").Formula = "=SUM(B5:B14)"

# 5) Match the cursor position Excel would leave behind after this edit.
$ws.Range("B15").Select()
